# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# The data now lives on the same sheet (columns AD:AF) instead of a
# separate sheet, as described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look exactly like the rest of the header row
# (bold, centered, thin-bordered). Copy the formatting from the last
# existing header cell (AC1) before writing the new header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-46) gets the same team record values.
$ws.Range("AD2:AD46").Value = 76
$ws.Range("AE2:AE46").Value = 85
$ws.Range("AF2:AF46").Value = 0
